$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("my_gt")
$ws.Range("A9").Value = "ide softmax 256 @ basis, -s 1"
$ws.Range("B9").Value = "60fps, 6epoch"
